$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-number-format on price cells that would otherwise be
# auto-coerced into numeric values by Excel, so they stay as literal text
# (matching the source data which stores prices as inline strings).
$textCells = @("D5","D6","D9","D10","D12","D15","D16","D21","D23","D24","D25","D26","D30","D32","D34","D35","D36","D38","D42","D43","D44","D46","D47","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "51.099.52"
$ws.Range("E2").Value = "  -0.12%  "
$ws.Range("D3").Value = "2.960.79"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "380.73"
$ws.Range("E5").Value = "  +1.43%  "
$ws.Range("D6").Value = "102.35"
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("E7").Value = "  +1.88%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").Value = "0.589"
$ws.Range("E9").Value = "  +0.98%  "
$ws.Range("D10").Value = "36.59"
$ws.Range("E10").Value = "  -0.44%  "
$ws.Range("E11").Value = "  -0.66%  "
$ws.Range("D12").Value = "0.0853"
$ws.Range("E12").Value = "  +2.04%  "
$ws.Range("D13").Value = "3.426.55"
$ws.Range("E13").Value = "  +0.73%  "
$ws.Range("E14").Value = "  +2.70%  "
$ws.Range("D15").Value = "7.75"
$ws.Range("E15").Value = "  +5.60%  "
$ws.Range("D16").Value = "12.10"
$ws.Range("E16").Value = "  +70.21%  "
$ws.Range("D17").Value = "2.964.60"
$ws.Range("E17").Value = "  +0.28%  "
$ws.Range("E18").Value = "  +3.31%  "
$ws.Range("D19").Value = "51.176.97"
$ws.Range("E19").Value = "  +0.18%  "
$ws.Range("E20").Value = "  -1.96%  "
$ws.Range("D21").Value = "12.45"
$ws.Range("E21").Value = "  -0.94%  "
$ws.Range("D22").Value = "0.0₃0963"
$ws.Range("E22").Value = "  +1.03%  "
$ws.Range("D23").Value = "3.34"
$ws.Range("E23").Value = "  +16.37%  "
$ws.Range("D24").Value = "269.38"
$ws.Range("E24").Value = "  +2.39%  "
$ws.Range("D25").Value = "69.80"
$ws.Range("E25").Value = "  +2.32%  "
$ws.Range("D26").Value = "7.92"
$ws.Range("E26").Value = "  -2.66%  "
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("E28").Value = "  -0.95%  "
$ws.Range("E29").Value = "  +0.90%  "
$ws.Range("D30").Value = "7.07"
$ws.Range("E30").Value = "  -9.81%  "
$ws.Range("D32").Value = "10.47"
$ws.Range("E32").Value = "  +6.15%  "
$ws.Range("E33").Value = "  +5.55%  "
$ws.Range("D34").Value = "51.38"
$ws.Range("E34").Value = "  +0.81%  "
$ws.Range("D35").Value = "34.38"
$ws.Range("E35").Value = "  +0.35%  "
$ws.Range("D36").Value = "0.0437"
$ws.Range("E36").Value = "  -4.17%  "
$ws.Range("E37").Value = "  -0.03%  "
$ws.Range("D38").Value = "3.32"
$ws.Range("E38").Value = "  +11.37%  "
$ws.Range("E39").Value = "  +1.95%  "
$ws.Range("E40").Value = "  +1.80%  "
$ws.Range("E41").Value = "  +3.53%  "
$ws.Range("D42").Value = "2.51"
$ws.Range("E42").Value = "  -2.15%  "
$ws.Range("D43").Value = "124.68"
$ws.Range("E43").Value = "  +2.36%  "
$ws.Range("D44").Value = "21.70"
$ws.Range("E44").Value = "  +2.95%  "
$ws.Range("E45").Value = "  +10.40%  "
$ws.Range("D46").Value = "0.277"
$ws.Range("E46").Value = "  +1.89%  "
$ws.Range("B47").Value = "WEMIXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Value = "2.02"
$ws.Range("E47").Value = "  -1.72%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "2.064.08"
$ws.Range("E48").Value = "  +3.42%  "
$ws.Range("E49").Value = "  +1.99%  "
$ws.Range("D50").Value = "0.0321"
$ws.Range("E50").Value = "  -8.75%  "
$ws.Range("D51").Value = "5.39"
$ws.Range("E51").Value = "  +7.19%  "

# Restore General number format on the cells we forced to text so the
# cell styling matches a normal (unformatted) text cell.
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "General"
}